$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$statusText = "Handback transform failed"
$zhError = "Handback file name: wsdfrqur.hbf is different with handoff file name: d1526d08-1c68-4559-ba5d-bfc544d9e3f9.74d059fe64636fa53de4700abc357c0021a2c8a2.zh-cn."
$deError = "Handback file name: wsdfrqur.hbf is different with handoff file name: d1526d08-1c68-4559-ba5d-bfc544d9e3f9.74d059fe64636fa53de4700abc357c0021a2c8a2.de-de."

# Update status text on the Overview summary sheet (row for the d1526d08... handoff)
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

# zh-cn sheet: status + error detail for the failed handback
$wsZhCn.Range("C3").Value = $statusText
$wsZhCn.Range("L3").Value = $zhError

# de-de sheet: status + error detail for the failed handback
$wsDeDe.Range("C3").Value = $statusText
$wsDeDe.Range("L3").Value = $deError
